# [Kadastro App] Yeni kayit eklendi: 2904
# Appends the new record (record no 2904) to the bottom of both the master
# "Kayitlar" list and the per-birim "Erdemli" list, mirroring how the app
# logs a newly created Erdemli kaydi.

$wb = $excel.ActiveWorkbook

$recordNo   = "2904"
$tarih      = "2025-09-08"
$birim      = "Erdemli"
$parselSayisi = ""
$is         = "3B"
$personellerKayitlar = "SEVİL SARAÇER (Tekniker), EMİNE ALANLI KIRCILI (K.Mühendisi)"
$personellerErdemli  = "SEVİL SARAÇER (Tekniker), EMİNE ALANLI KIRCILI (K.Mühendisi)"

function Add-Kayit($ws, $row, $personeller) {
    # Columns A (Kayit No), B (Tarih) and D (Parsel Sayisi) hold values that
    # look like numbers/dates ("2904", "2025-09-08", blank) - force them to
    # Text so they are stored the same way as the rest of the (text-only)
    # column instead of being auto-converted to a real number/date.
    $numericLookingCells = $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 2))
    $numericLookingCells.NumberFormat = "@"
    $ws.Cells.Item($row, 4).NumberFormat = "@"

    $ws.Cells.Item($row, 1).Value = $recordNo
    $ws.Cells.Item($row, 2).Value = $tarih
    $ws.Cells.Item($row, 3).Value = $birim
    $ws.Cells.Item($row, 4).Value = $parselSayisi
    $ws.Cells.Item($row, 5).Value = $is
    $ws.Cells.Item($row, 6).Value = $personeller
}

# "Kayitlar" (sheet1) master list -> new row 11
$wsKayitlar = $wb.Worksheets.Item("Kayitlar")
Add-Kayit $wsKayitlar 11 $personellerKayitlar

# "Erdemli" (sheet8) birim list -> new row 10
$wsErdemli = $wb.Worksheets.Item("Erdemli")
Add-Kayit $wsErdemli 10 $personellerErdemli
